$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new column for "d=6" right before the current column G
#    (which holds "d=7"). This shifts d=7: G->H and d=10: H->I.
# ------------------------------------------------------------------
$ws.Columns("G").Insert()

# ------------------------------------------------------------------
# 2) Insert the new data rows for the "d=6" results. They need to be
#    inserted right before the "ARMA_I(x,7,0)" row of each AR block,
#    which (after the column insert above) sit at rows 20, 35 and 49.
#    We insert from the bottom-most block upward so the row numbers
#    used for the still-untouched blocks stay valid.
# ------------------------------------------------------------------
$ws.Rows("49:50").Insert()   # room for ARMA_I(2,6,0) / ARMA_I(2,6,2)
$ws.Rows("35:36").Insert()   # room for ARMA_I(1,6,0) / ARMA_I(1,6,1)
$ws.Rows("20:22").Insert()   # room for ARMA_I(0,6,0..2)

# ------------------------------------------------------------------
# 3) New header cell for the inserted column
# ------------------------------------------------------------------
$ws.Range("G1").Value = "d=6"

# ------------------------------------------------------------------
# 4) Fill in the new rows (label in column A, value in column G) and
#    give column A the same look (bold, centered, boxed) as the rest
#    of the table by copying the format from the row above.
# ------------------------------------------------------------------
function Set-NewRow($row, $label, $value) {
    $ws.Range("A$row").Value = $label
    $ws.Range("G$row").Value = $value

    $srcRow = $row - 1
    $ws.Range("A$srcRow").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
}

Set-NewRow 20 "ARMA_I(0,6,0)" 98.00057275652443
Set-NewRow 21 "ARMA_I(0,6,1)" 98.03625580966863
Set-NewRow 22 "ARMA_I(0,6,2)" 97.95451476935884

Set-NewRow 38 "ARMA_I(1,6,0)" 97.97946172497079
Set-NewRow 39 "ARMA_I(1,6,1)" 97.92550230539541

Set-NewRow 54 "ARMA_I(2,6,0)" 97.95374976245849
Set-NewRow 55 "ARMA_I(2,6,2)" 97.89261788502752

$excel.Application.CutCopyMode = $false

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
